$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMM")
$ws.Activate() | Out-Null

# Resource Management row (row 11) — rewording from "accounts/licenses"
# specific phrasing to a more general "resources" phrasing (licenses,
# accounts, domains, recurring expenses, etc.)
$ws.Range("C11").Value = "Resources like licenses, accounts, or domains are only tracked upon reminder of expiration or renewal needs; ownership is dispersed across multiple people"
$ws.Range("D11").Value = "One person tracks resources; knowledge not available to the entire Red Team"
$ws.Range("E11").Value = "Resources are centrally tracked, understood, and reviewed as needed by the Red Team; Red Team account passwords are secured"
$ws.Range("F11").Value = "Recurring expenses or other resources are reviewed quarterly for need or expiration"

# The new wording is longer, so the wrapped-text row grows taller once
# re-flowed (mirrors Excel's automatic row re-height on content change).
$ws.Rows("11:11").RowHeight = 75.75

# Leave the view scrolled/selected/zoomed where the edit was made.
$ws.Range("G11").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
